# The deck's Slide Master has 6 custom layouts. In the "before" state the
# "Blank" layout sits first; the commit reverts to an earlier arrangement
# where "Blank" is the last layout and everything else shifts up by one
# slot (Custom Layout, Primera Lectura, Salmo, Segunda Lectura,
# 2_Primera Lectura, Blank).
#
# Reproduce that by moving the "Blank" custom layout from position 1 to
# the end of the Slide Master's CustomLayouts collection.

$p = $ppt.ActivePresentation
$m = $p.SlideMaster

$blank = $m.CustomLayouts.Item(1)
$blank.MoveTo($m.CustomLayouts.Count)
